$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column G (EOT_response) to fit its contents
$ws.Columns("G:G").ColumnWidth = 13.6

# Zoom the active window in (matches the workbook's saved view state)
$excel.ActiveWindow.Zoom = 183

# Append 10 new data rows (rows 12-21)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 0.0111
$ws.Range("D12").Value = 0.215
$ws.Range("E12").Value = 1.02
$ws.Range("F12").Value = 0.9757
$ws.Range("G12").Value = "Non_Responder"
$ws.Range("H12").Value = "Training"

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = 0.0165
$ws.Range("D13").Value = 0.205
$ws.Range("E13").Value = 0.99
$ws.Range("F13").Value = 0.9763
$ws.Range("G13").Value = "Responder"

$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 0.01
$ws.Range("D14").Value = 0.18
$ws.Range("E14").Value = 0.97
$ws.Range("F14").Value = 0.9765
$ws.Range("G14").Value = "Responder"
$ws.Range("H14").Value = "Training"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 0.018
$ws.Range("D15").Value = 0.22
$ws.Range("E15").Value = 1.25
$ws.Range("F15").Value = 0.9759
$ws.Range("G15").Value = "Non_Responder"
$ws.Range("H15").Value = "Training"

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = 3
$ws.Range("C16").Value = 0.021
$ws.Range("D16").Value = 0.195
$ws.Range("E16").Value = 0.89
$ws.Range("F16").Value = 0.976
$ws.Range("G16").Value = "Responder"

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = 2
$ws.Range("C17").Value = 0.03
$ws.Range("D17").Value = 0.24
$ws.Range("E17").Value = 1.1
$ws.Range("F17").Value = 0.9755
$ws.Range("G17").Value = "Responder"

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 0.0089
$ws.Range("D18").Value = 0.2001
$ws.Range("E18").Value = 0.95
$ws.Range("F18").Value = 0.9761
$ws.Range("G18").Value = "Responder"
$ws.Range("H18").Value = "Validation"

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = 0.02
$ws.Range("D19").Value = 0.23
$ws.Range("E19").Value = 1.3
$ws.Range("F19").Value = 0.9762
$ws.Range("G19").Value = "Non_Responder"

$ws.Range("A20").Value = 19
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 0.015
$ws.Range("D20").Value = 0.2105
$ws.Range("E20").Value = 0.76
$ws.Range("F20").Value = 0.9758
$ws.Range("G20").Value = "Non_Responder"
$ws.Range("H20").Value = "Training"

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 0.012345
$ws.Range("D21").Value = 0.19
$ws.Range("E21").Value = 1.2
$ws.Range("F21").Value = 0.976
$ws.Range("G21").Value = "Responder"
$ws.Range("H21").Value = "Validation"
